# Applies the crypto price/volume update described by the commit diff.
# Numeric-looking values (e.g. "251.33") must stay TEXT cells (matching the
# source workbook, which stores every data cell as an inline string) --
# Excel's COM Range.Value setter auto-coerces parseable numerics, so for
# those we momentarily force a Text number format, assign, then clear the
# format again so the cell keeps its original (default) style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.147.23"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.012.76"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.642"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.53"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +17.31%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.07"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.41%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.929"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.89"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").Value = "2.309.09"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +15.38%  "
$ws.Range("D18").Value = "1.985.73"
$ws.Range("E18").Value = "  -3.22%  "
$ws.Range("D19").Value = "36.061.55"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").Value = "0.0₃0855"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +23.78%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.120"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.68%  "
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +24.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.47"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +13.15%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.73"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +17.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +22.03%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.88"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.29%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "93.99"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.81"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.29%  "
$ws.Range("D48").Value = "1.422.16"
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.91"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.95%  "
